$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.118.58'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '3.164.67'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D8').Value = '3.162.04'
$ws.Range('E8').Value = '  -1.57%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.544'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.36%  '
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('E11').Value = '  -8.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.516'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.77%  '
$ws.Range('E13').Value = '  -1.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.21'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').Value = '3.687.37'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').Value = '66.146.26'
$ws.Range('E16').Value = '  +0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.39'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('D18').Value = '3.172.35'
$ws.Range('E18').Value = '  -2.10%  '
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '509.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.725'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.07'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.36'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.12'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.88'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.76%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  -2.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.49'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.46%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '499.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0876'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0418'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.127'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.73'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.98%  '
$ws.Range('D42').Value = '0.0₃0679'
$ws.Range('E42').Value = '  +6.36%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.295'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.79'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.15%  '
$ws.Range('E45').Value = '  -1.57%  '
$ws.Range('D46').Value = '2.819.96'
$ws.Range('E46').Value = '  -4.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.36'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.116'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.46%  '
